# Fixed update to excel issue
# - Rename "Requested quantity" headers to dataset-specific names
# - Add a new "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper)

$wb = $excel.ActiveWorkbook

# --- Update existing sheet headers -----------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---------
# Duplicate the last sheet so the new sheet inherits the same sheet-level
# properties (outline settings, page setup, margins, header/date styles)
# used throughout the workbook, then wipe its contents and fill in the
# forecast data.
$wsMonthly.Copy($null, $wsMonthly)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "PO Forecast"
$ws.Cells.ClearContents()

$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Extend the bold/centered/bordered header formatting (already present on
# A1:B1 from the copied sheet) across the new C1:D1 header cells.
$ws.Range("A1:B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# Extend the date-formatted style (already present on A2 from the copied
# sheet) down through the rest of column A's data rows.
$ws.Range("A2").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)

$ws.Range("A2").Value = 45543.99999999999
$ws.Range("B2").Value = 29
$ws.Range("C2").Value = 23.26446051973556
$ws.Range("D2").Value = 33.93049518325795
$ws.Range("A3").Value = 45564.99999999999
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 16.16631454903294
$ws.Range("D3").Value = 27.63128486917476
$ws.Range("A4").Value = 45592.99999999999
$ws.Range("B4").Value = 13
$ws.Range("C4").Value = 7.879663218008544
$ws.Range("D4").Value = 18.91497931162268
$ws.Range("A5").Value = 45599.99999999999
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 5.650455342774682
$ws.Range("D5").Value = 16.1247752276285
$ws.Range("A6").Value = 45606.99999999999
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 3.368212828955232
$ws.Range("D6").Value = 14.5430187468413
$ws.Range("A7").Value = 45613.99999999999
$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 1.488405267904641
$ws.Range("D7").Value = 12.35855054009215
$ws.Range("A8").Value = 45620.99999999999
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = -0.5043379223704224
$ws.Range("D8").Value = 10.25158861093846
$ws.Range("A9").Value = 45627.99999999999
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = -2.713339376636737
$ws.Range("D9").Value = 7.932088090523981
$ws.Range("A10").Value = 45634.99999999999
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = -4.809279118449519
$ws.Range("D10").Value = 5.715020662342249
$ws.Range("A11").Value = 45641.99999999999
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = -7.448554336688312
$ws.Range("D11").Value = 3.838846743759725
$ws.Range("A12").Value = 45648.99999999999
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = -8.851503552373714
$ws.Range("D12").Value = 1.711470134622908

$ws.Range("A1").Select()

Write-Host "PO Forecast sheet added and headers updated."
